$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subgroups")

$ws.Range("A14").Value = '{''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A16").Value = '{''Hobby'': ''1'', ''DevType'': ''2''}'
$ws.Range("A17").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A19").Value = '{''Dependents'': ''2'', ''Hobby'': ''1''}'
$ws.Range("A20").Value = '{''Age'': ''3'', ''Hobby'': ''1''}'
$ws.Range("A23").Value = '{''UndergradMajor'': ''2'', ''Student'': ''1''}'
$ws.Range("A24").Value = '{''Student'': ''1'', ''DevType'': ''2''}'
$ws.Range("A27").Value = '{''Dependents'': ''2'', ''Student'': ''1''}'
$ws.Range("A28").Value = '{''Age'': ''3'', ''Student'': ''1''}'
$ws.Range("A29").Value = '{''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A30").Value = '{''GINI'': ''2'', ''Student'': ''1''}'
$ws.Range("A31").Value = '{''FormalEducation'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A32").Value = '{''UndergradMajor'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A33").Value = '{''UndergradMajor'': ''2'', ''RaceEthnicity'': ''1''}'
$ws.Range("A38").Value = '{''HDI'': ''1'', ''DevType'': ''2''}'
$ws.Range("A40").Value = '{''Dependents'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A41").Value = '{''Age'': ''3'', ''SexualOrientation'': ''1''}'
$ws.Range("A42").Value = '{''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A43").Value = '{''GINI'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A46").Value = '{''Dependents'': ''2'', ''HDI'': ''1''}'
$ws.Range("A47").Value = '{''HDI'': ''1'', ''GDP'': ''1''}'
$ws.Range("A48").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''DevType'': ''2''}'
$ws.Range("A49").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A50").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A51").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A52").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A53").Value = '{''UndergradMajor'': ''2'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A54").Value = '{''Hobby'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A55").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A56").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A57").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A58").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''RaceEthnicity'': ''1''}'
$ws.Range("A60").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''HDI'': ''1''}'
$ws.Range("A61").Value = '{''UndergradMajor'': ''2'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A62").Value = '{''Student'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A64").Value = '{''Dependents'': ''2'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A65").Value = '{''HDI'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A66").Value = '{''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''Student'': ''1''}'
$ws.Range("A67").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A68").Value = '{''Dependents'': ''2'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A69").Value = '{''UndergradMajor'': ''2'', ''RaceEthnicity'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A70").Value = '{''UndergradMajor'': ''2'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A72").Value = '{''HDI'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A73").Value = '{''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A74").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A75").Value = '{''Dependents'': ''2'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A76").Value = '{''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''HDI'': ''1''}'
$ws.Range("A77").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A78").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A79").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A80").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A81").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A82").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A83").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
